# Add a new translation entry for "ScaledText" just above the existing
# "EditText" row (old row 96), shifting all subsequent rows down by one.
#
# Column layout: A = ID, B = ENGLISH, C = SPANISH

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 96; everything currently at/after row 96
# (EditText, CurrentMode, ... EventsDoorStateTooltip) moves down to row 97+.
$ws.Rows.Item(96).Insert()

# Fill the new row with the "Scaled Text" translation triplet.
$ws.Range("A96").Value = "ScaledText"
$ws.Range("B96").Value = "Scaled Text"
$ws.Range("C96").Value = "Texto Escalado"

# Match the author's final selection/viewport position recorded in the diff.
[void]$ws.Range("D96").Select()
